$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1) Update the "time_taken" timestamps (column F, rows 2-18) on the "data"
#    sheet.
# ---------------------------------------------------------------------------
$dataWs = $wb.Worksheets.Item("data")

$timestamps = @{
    2  = "2021-10-05 14:35:34.721564"
    3  = "2021-10-05 14:35:34.721572"
    4  = "2021-10-05 14:35:34.721575"
    5  = "2021-10-05 14:35:34.721578"
    6  = "2021-10-05 14:35:34.721581"
    7  = "2021-10-05 14:35:34.721584"
    8  = "2021-10-05 14:35:34.721586"
    9  = "2021-10-05 14:35:34.721589"
    10 = "2021-10-05 14:35:34.721592"
    11 = "2021-10-05 14:35:34.721595"
    12 = "2021-10-05 14:35:34.721597"
    13 = "2021-10-05 14:35:34.721600"
    14 = "2021-10-05 14:35:34.721602"
    15 = "2021-10-05 14:35:34.721605"
    16 = "2021-10-05 14:35:34.721607"
    17 = "2021-10-05 14:35:34.721610"
    18 = "2021-10-05 14:35:34.721612"
}

foreach ($row in $timestamps.Keys) {
    $dataWs.Cells.Item($row, 6).Value = $timestamps[$row]
}

# ---------------------------------------------------------------------------
# 2) Add a new "metadata" worksheet right after "data".
#
#    We duplicate the "data" sheet (which carries over the sheet-level
#    formatting -- outline/page-setup properties, margins, etc.) and then
#    wipe its contents, rather than starting from a completely blank sheet.
# ---------------------------------------------------------------------------
$dataWs.Copy([System.Reflection.Missing]::Value, $dataWs)
$newWs = $wb.Worksheets.Item(2)
$newWs.Name = "metadata"
$newWs.Cells.Clear()

# Header row
$newWs.Cells.Item(1, 2).Value = "data_name"
$newWs.Cells.Item(1, 3).Value = "data_id"
$newWs.Cells.Item(1, 4).Value = "data_version"
$newWs.Cells.Item(1, 5).Value = "data_version_created"
$newWs.Cells.Item(1, 6).Value = "panel_query_time"
$newWs.Cells.Item(1, 7).Value = "panel_get_request"

# Re-use the exact header style/border/alignment already defined in the
# workbook (same look as the "data" sheet header row) by copying formats
# across instead of re-building a style from scratch.
$dataWs.Range("B1:F1").Copy()
$newWs.Range("B1:G1").PasteSpecial(-4122)  # xlPasteFormats

# Data row
$newWs.Cells.Item(2, 1).Value = 0
$dataWs.Range("A2").Copy()
$newWs.Range("A2").PasteSpecial(-4122)     # xlPasteFormats

$newWs.Cells.Item(2, 2).Value = "Renal Hypertension and Disorders of Aldosterone Metabolism"
$newWs.Cells.Item(2, 3).Value = 190

# Force this cell to remain text ("1.4") instead of being auto-converted to
# the number 1.4.
$newWs.Cells.Item(2, 4).NumberFormat = "@"
$newWs.Cells.Item(2, 4).Value = "1.4"

$newWs.Cells.Item(2, 5).Value = "2021-08-05T01:49:53.841164Z"
$newWs.Cells.Item(2, 6).Value = "2021-10-05 14:35:34.717911"
$newWs.Cells.Item(2, 7).Value = "https://panelapp.agha.umccr.org/api/v1/panels/190/?format=json"

$newWs.Range("A1").Select() | Out-Null
$dataWs.Activate() | Out-Null
